# Weekly fruit/vegetable price update: insert one new daily-price record
# for "Naranja" (orange), variety "Valencia", quality "Primera" at
# Terminal Hortofrutícola Agro Chillán, shifting the existing rows
# 536-575 down to 537-576.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 536 (pushes old row 536..575 -> 537..576)
$ws.Rows.Item(536).Insert()

# Populate the new row with the new record. Most columns (market, region,
# product taxonomy, unit, origin, etc.) repeat the values already used for
# this Valencia / Primera series; only the date, volume and the three
# price columns (+ derived $/Kg) change.
$ws.Cells.Item(536, 1).Value = 7
$ws.Cells.Item(536, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(536, 3).Value = "Ñuble"
$ws.Cells.Item(536, 4).Value = 45013
$ws.Cells.Item(536, 5).Value = 16
$ws.Cells.Item(536, 6).Value = "Fruta"
$ws.Cells.Item(536, 7).Value = 100102
$ws.Cells.Item(536, 8).Value = "Cítricos"
$ws.Cells.Item(536, 9).Value = 100102005
$ws.Cells.Item(536, 10).Value = "Naranja"
$ws.Cells.Item(536, 11).Value = "Valencia"
$ws.Cells.Item(536, 12).Value = "Primera"
$ws.Cells.Item(536, 13).Value = 80
$ws.Cells.Item(536, 14).Value = 13000
$ws.Cells.Item(536, 15).Value = 13000
$ws.Cells.Item(536, 16).Value = 13000
$ws.Cells.Item(536, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(536, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(536, 19).Value = 867
$ws.Cells.Item(536, 20).Value = 15

# Keep the same date-time number format on column D (row 536) as the rest
# of the "Fecha" column uses.
$ws.Cells.Item(536, 4).NumberFormat = $ws.Cells.Item(537, 4).NumberFormat
